$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnes")

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "bernardcac"
